# Insert a new data row above row 250 (shifts existing rows 250:360 down to 251:361)
# and populate the new row 250 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(250).Insert()

$ws.Range("A250").Value = 8
$ws.Range("B250").Value = "Terminal La Palmera de La Serena"
$ws.Range("C250").Value = "Coquimbo"
$ws.Range("D250").Value = 44726
$ws.Range("D250").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E250").Value = 4
$ws.Range("F250").Value = 100114013
$ws.Range("G250").Value = "Zanahoria"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 480
$ws.Range("K250").Value = 6000
$ws.Range("L250").Value = 7000
$ws.Range("M250").Value = 6500
$ws.Range("N250").Value = "`$/saco 20 kilos"
$ws.Range("O250").Value = "Provincia del Elquí"
$ws.Range("P250").Value = 325
$ws.Range("Q250").Value = 20
$ws.Range("R250").Value = "Hortaliza"
